# The document ends with a paragraph that holds both the final sentence
# of bullet text and the (hidden) _GoBack bookmark, all in one <w:p>:
#
#   ...va de esa manera.[bookmarkStart _GoBack][bookmarkEnd _GoBack]<pilcrow>
#
# The edit splits that into two paragraphs: the sentence keeps its own
# paragraph, and a new paragraph (no list style/numbering, just bold +
# es-ES run formatting on its mark) holding only the _GoBack bookmark
# follows it.

$d = $word.ActiveDocument

# Find the end of the final sentence so we split right after its
# trailing period, i.e. right before the bookmark that sits at the very
# end of the story. Use a short, punctuation-free needle so it matches
# regardless of any smart-quote characters elsewhere in the sentence.
# Find.Execute collapses/resizes the Range it was called on to the
# match, so re-use that same Range object to read back its new bounds
# (re-fetching $d.Content would just return the whole-story range again).
$searchRange = $d.Content
$needle = "va de esa manera."
$found = $searchRange.Find.Execute($needle, $true, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)
if (-not $found) {
  throw "could not locate the closing sentence to split on"
}
$sentenceEnd = $searchRange.End

# Splitting exactly at the bookmark's position leaves the bookmark glued
# to the first paragraph, so split one character earlier instead (right
# before the trailing "."); the bookmark/pilcrow then lands cleanly on
# the newly created paragraph, inheriting the source paragraph's bullet
# style/numbering and bold/es-ES run formatting.
$splitPoint = $sentenceEnd - 1
$d.Range($splitPoint, $splitPoint).InsertParagraphAfter()

# That new (now last) paragraph contains "." + the bookmark; replace its
# contents with just the bookmark, keeping the bold/es-ES paragraph-mark
# formatting but dropping the inherited bulleted-list style/numbering.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$newPara.Range.InsertXML($bookmarkXml)

# Restore the trailing period that was removed from the sentence
# paragraph when the split point was backed up by one character.
$textPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$periodPos = $textPara.Range.End - 1
$d.Range($periodPos, $periodPos).InsertAfter(".")
